# Daily auto push 2026-01-06 04:09 UTC
# A new measurement row (2026/01/06, 12:00, rank 18) lands between the existing
# 2026/01/06 09:00 entry and the 2026/12/29 block, so every row from the old
# 567 through 608 shifts down by one and the sheet grows by a single row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row at 567; Excel shifts rows 567:608 down to 568:609
# and extends the used range/dimension automatically.
$ws.Rows("567").Insert()

# Populate the freshly inserted row 567 with the new data point.
# The leading apostrophe forces Excel to store the date-looking text as a
# literal string instead of auto-converting it to a date serial number
# (matching every other cell in column A, which is stored as text).
$ws.Range("A567").Value = "'2026/01/06"
# Re-apply the default "Normal" style so no stray number-format style index
# is left attached to the cell (keeps formatting identical to its neighbours).
$ws.Range("A567").Style = "Normal"

$ws.Range("B567").Value = "火"
$ws.Range("C567").Value = 12
$ws.Range("D567").Value = 18
